$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# 1. Delete the "packages" sheet entirely (no longer used).
$wb.Worksheets("packages").Delete() | Out-Null

# 2. Reorder: move "attributes" to be the first sheet (before "entities").
#    Final order: attributes, entities, importperson, importcity
$wb.Worksheets("attributes").Move($wb.Worksheets("entities"))

# 3. Rename the import sheets to use underscores.
$wb.Worksheets("importperson").Name = "import_person"
$wb.Worksheets("importcity").Name = "import_city"

# 4. Fix values on the "attributes" sheet: the entity-name column (B) and the
#    refEntity column (E) still used the old "importperson"/"importcity" names.
$wsAttr = $wb.Worksheets("attributes")
$wsAttr.Range("B2").Value = "import_city"
$wsAttr.Range("B3").Value = "import_person"
$wsAttr.Range("B4").Value = "import_person"
$wsAttr.Range("B5").Value = "import_person"
$wsAttr.Range("B6").Value = "import_person"
$wsAttr.Range("B7").Value = "import_person"
$wsAttr.Range("B8").Value = "import_person"
$wsAttr.Range("E7").Value = "import_person"
$wsAttr.Range("E8").Value = "import_city"

# 5. Fix the "entities" sheet: column B held a stray/incorrect "package"/"test"
#    column that duplicated column C; drop it so the real description column
#    shifts into B, then fix up the renamed entity names in column A.
$wsEnt = $wb.Worksheets("entities")
$wsEnt.Columns("B").Delete() | Out-Null
$wsEnt.Range("A2").Value = "import_city"
$wsEnt.Range("A3").Value = "import_person"

# 6. Restore selections / active sheet to match the reverted state.
$wsAttr.Range("A6:XFD6").Select() | Out-Null

$wsEnt.Range("A4:XFD4").Select() | Out-Null

$wsPerson = $wb.Worksheets("import_person")
$wsPerson.Range("F7").Select() | Out-Null
$wsPerson.Activate()
